# feat: OPTICS & eigendocs vs. 32x32
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Task (column B) text updates on existing rows ---
$ws.Range("B37").Value = "fuzzy full-text search, DB: Elasticsearch"

# --- Rows 39-41 gain Date (A) + Task (B) entries (column D stays as-is) ---
# Copy the Date/Task formatting (borders, number format, style) from row 38
# without touching column D's existing content.
$ws.Range("A38:B38").Copy()
$ws.Range("A39:B39").PasteSpecial(-4122)
$ws.Range("A40:B40").PasteSpecial(-4122)
$ws.Range("A41:B41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A39").Value = 45186
$ws.Range("B39").Value = "BA: Clustering"

$ws.Range("A40").Value = 45187
$ws.Range("B40").Value = "OPTICS: reachability plot + clustering"

$ws.Range("A41").Value = 45188
$ws.Range("B41").Value = "OPTICS: on 32x32 and on PCA version (PCA/eigendoc is better), PCA on 2^2, 5^2, 14^2 -> 2^2 is best, eigendocs code und BA"

$ws.Rows.Item(39).RowHeight = 18
$ws.Rows.Item(40).RowHeight = 18
$ws.Rows.Item(41).RowHeight = 52

# --- TODO (column D) text updates / consolidation ---
$ws.Range("D43").Value = "reduce tfidf vocab "
$ws.Range("D46").Value = "OPTICS: threshold hierarchical clustering, reduce threshold and observe new documents -> schwer"
$ws.Range("D47").Value = "database on cluster"

# --- Remove the now-superseded D-only rows (their notes were folded into the
#     consolidated entries above). These rows are not shifted up - they are
#     simply wiped so they disappear from the saved sheet, leaving rows
#     43, 46 and 47 at their original row numbers. ---
$ws.Range("D42").Clear()
$ws.Rows.Item(42).EntireRow.AutoFit()
$ws.Range("D44").Clear()
$ws.Rows.Item(44).EntireRow.AutoFit()
$ws.Range("D45").Clear()
$ws.Rows.Item(45).EntireRow.AutoFit()

# --- Move the cursor to match the author's final selection ---
$ws.Range("D46").Select()
